{"js": "// Remove the \"Inspect code for aliasing violations.\" to-do list item.\n//\n// Locate the paragraph by its text (robust to its position shifting) and\n// delete the whole paragraph (including its paragraph mark / list bullet),\n// leaving the surrounding list items untouched.\nconst body = context.document.body;\n\nconst results = body.search(\"Inspect code for aliasing violations.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the paragraph to remove.\");\n}\n\n// Delete every matching paragraph (normally just the one).\nfor (const hit of results.items) {\n  const para = hit.paragraphs.getFirst();\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Inspect code for aliasing violations.\" to-do list item.\n#\n# Locate the paragraph with Find (robust to its position shifting), expand\n# the hit to the full paragraph (wdParagraph = 4, includes the paragraph\n# mark / list bullet) and delete it, leaving the surrounding list items\n# untouched.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$found = $range.Find.Execute(\"Inspect code for aliasing violations.\")\n\nif ($found) {\n    $range.Expand(4)  # wdParagraph\n    $range.Delete()\n}\n"}
